$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the weekly progress report fields
$ws.Range("B2").Value = 7

$ws.Range("B3").Value = 45977   # 11/16/2025 (date serial number)
$ws.Range("B3").NumberFormat = "mm-dd-yy"

$ws.Range("B4").Value = 45984   # 11/23/2025 (date serial number)
# Reuse B3's date style for B4 instead of creating a duplicate number format entry
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Range("B5").Value = "Lại Việt Anh"
$ws.Range("B6").Value = 179066
$ws.Range("B7").Value = "Xây dựng mô hình giám sát và điều khiển nhà thông minh"
$ws.Range("B8").Value = "thiết kế hệ thống"
$ws.Range("B9").Value = "thiết kế hệ thống"

# Rows 8 and 9 no longer need the tall wrapped-text height now that the
# content fits on a single line.
$ws.Rows.Item(8).RowHeight = 15.6
$ws.Rows.Item(9).RowHeight = 15.6

$ws.Range("I9").Select()
